$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the display-name values for the two questions (text only change,
# as part of attaching images to the questions).
$ws.Range("A3").Value = "haha"
$ws.Range("A2").Value = "hohoh"

# Update the last active selection to match the edited workbook.
$ws.Range("A10").Select()
